$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay stored as text (matching source data),
# mirroring how these cells already held text before the edit.
$textRows = @(5,6,11,12,13,14,16,20,21,22,23,24,25,26,27,28,30,31,32,33,36,37,38,40,41,42,44,45,50,51)
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = "66.806.52"
$ws.Cells.Item(2, 5).Value = "  +0.70%  "
$ws.Cells.Item(3, 4).Value = "3.268.57"
$ws.Cells.Item(3, 5).Value = "  -2.42%  "
$ws.Cells.Item(4, 5).Value = "  +0.62%  "
$ws.Cells.Item(5, 4).Value = "575.35"
$ws.Cells.Item(5, 5).Value = "  -1.89%  "
$ws.Cells.Item(6, 4).Value = "171.02"
$ws.Cells.Item(6, 5).Value = "  -8.02%  "
$ws.Cells.Item(7, 5).Value = "  +0.21%  "
$ws.Cells.Item(8, 5).Value = "  -0.09%  "
$ws.Cells.Item(9, 4).Value = "3.261.23"
$ws.Cells.Item(9, 5).Value = "  -2.51%  "
$ws.Cells.Item(10, 5).Value = "  -5.98%  "
$ws.Cells.Item(11, 4).Value = "0.567"
$ws.Cells.Item(11, 5).Value = "  -2.87%  "
$ws.Cells.Item(12, 4).Value = "44.67"
$ws.Cells.Item(12, 5).Value = "  -5.21%  "
$ws.Cells.Item(13, 4).Value = "0.0000268"
$ws.Cells.Item(13, 5).Value = "  -0.58%  "
$ws.Cells.Item(14, 4).Value = "685.68"
$ws.Cells.Item(14, 5).Value = "  +2.45%  "
$ws.Cells.Item(15, 4).Value = "3.801.27"
$ws.Cells.Item(15, 5).Value = "  -1.59%  "
$ws.Cells.Item(16, 4).Value = "8.18"
$ws.Cells.Item(16, 5).Value = "  -4.09%  "
$ws.Cells.Item(17, 4).Value = "66.896.62"
$ws.Cells.Item(17, 5).Value = "  +0.78%  "
$ws.Cells.Item(18, 5).Value = "  +0.49%  "
$ws.Cells.Item(19, 4).Value = "3.279.26"
$ws.Cells.Item(19, 5).Value = "  -1.56%  "
$ws.Cells.Item(20, 4).Value = "17.09"
$ws.Cells.Item(20, 5).Value = "  -4.64%  "
$ws.Cells.Item(21, 4).Value = "10.59"
$ws.Cells.Item(21, 5).Value = "  -4.83%  "
$ws.Cells.Item(22, 4).Value = "0.877"
$ws.Cells.Item(22, 5).Value = "  -2.57%  "
$ws.Cells.Item(23, 4).Value = "16.78"
$ws.Cells.Item(23, 5).Value = "  -5.92%  "
$ws.Cells.Item(24, 4).Value = "5.16"
$ws.Cells.Item(24, 5).Value = "  +2.28%  "
$ws.Cells.Item(25, 4).Value = "97.92"
$ws.Cells.Item(25, 5).Value = "  -4.18%  "
$ws.Cells.Item(26, 4).Value = "3.81"
$ws.Cells.Item(26, 5).Value = "  -4.55%  "
$ws.Cells.Item(27, 4).Value = "2.61"
$ws.Cells.Item(27, 5).Value = "  -6.29%  "
$ws.Cells.Item(28, 4).Value = "33.16"
$ws.Cells.Item(28, 5).Value = "  +2.70%  "
$ws.Cells.Item(29, 5).Value = "  -4.65%  "
$ws.Cells.Item(30, 4).Value = "8.24"
$ws.Cells.Item(30, 5).Value = "  -3.41%  "
$ws.Cells.Item(31, 4).Value = "6.55"
$ws.Cells.Item(31, 5).Value = "  -4.71%  "
$ws.Cells.Item(32, 4).Value = "576.11"
$ws.Cells.Item(32, 5).Value = "  -5.84%  "
$ws.Cells.Item(33, 4).Value = "10.76"
$ws.Cells.Item(33, 5).Value = "  -3.34%  "
$ws.Cells.Item(34, 4).Value = "3.825.56"
$ws.Cells.Item(34, 5).Value = "  -0.78%  "
$ws.Cells.Item(35, 5).Value = "  -0.13%  "
$ws.Cells.Item(36, 4).Value = "0.101"
$ws.Cells.Item(36, 5).Value = "  -3.94%  "
$ws.Cells.Item(37, 4).Value = "55.29"
$ws.Cells.Item(37, 5).Value = "  -1.76%  "
$ws.Cells.Item(38, 4).Value = "3.26"
$ws.Cells.Item(38, 5).Value = "  -16.78%  "
$ws.Cells.Item(39, 5).Value = "  -0.92%  "
$ws.Cells.Item(40, 4).Value = "3.36"
$ws.Cells.Item(40, 5).Value = "  -1.14%  "
$ws.Cells.Item(41, 4).Value = "2.54"
$ws.Cells.Item(41, 5).Value = "  -5.48%  "
$ws.Cells.Item(42, 4).Value = "31.23"
$ws.Cells.Item(42, 5).Value = "  -5.60%  "
$ws.Cells.Item(43, 4).Value = "0.0₃0655"
$ws.Cells.Item(43, 5).Value = "  -6.93%  "
$ws.Cells.Item(44, 4).Value = "0.321"
$ws.Cells.Item(44, 5).Value = "  -5.15%  "
$ws.Cells.Item(45, 4).Value = "2.94"
$ws.Cells.Item(45, 5).Value = "  -8.24%  "
$ws.Cells.Item(46, 5).Value = "  -4.54%  "
$ws.Cells.Item(47, 5).Value = "  +0.05%  "
$ws.Cells.Item(48, 5).Value = "  -1.75%  "
$ws.Cells.Item(49, 5).Value = "  -1.01%  "
$ws.Cells.Item(50, 4).Value = "1.34"
$ws.Cells.Item(50, 5).Value = "  +3.29%  "
$ws.Cells.Item(51, 4).Value = "128.36"
$ws.Cells.Item(51, 5).Value = "  -0.79%  "
